# backup_projetos.xlsx edit:
#  - "Artista Cida Show" rows (15-27): mark tasks as Feito and fill in the
#    mined Mapa Cultural agent info (mapa URL, foto perfil, nome no mapa).
#  - "Artista Elisangela Monteiro" rows (380-392): update the proponente's
#    foto perfil URL and proper-case her name on the mapa.
#  - Remove the leftover "AA TESTE" scratch rows (406-409).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Artista Cida Show: rows 15 through 27 -------------------------------
$cidaMapa  = "https://www.mapacultural.pe.gov.br/agente/593/"
$cidaFoto  = "https://www.mapacultural.pe.gov.br/files/agent/593/file/1151869/blob.-2-5a26c036add54f2b5171661a842d970e.png"
$cidaNome  = "SHIRLEYDE ALBUQUERQUE"

for ($r = 15; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "Feito"
    $ws.Cells.Item($r, 5).Value = $cidaMapa
    $ws.Cells.Item($r, 6).Value = $cidaFoto
    $ws.Cells.Item($r, 7).Value = $cidaNome
}

# --- Artista Elisangela Monteiro: rows 380 through 392 --------------------
$elisFoto = "https://www.mapacultural.pe.gov.br/files/agent/16301/file/1170727/blob.-2-e50ddac85ec7945ec2518f56a4d1b778.png"
$elisNome = "Elisangela Monteiro De Melo Costa"

for ($r = 380; $r -le 392; $r++) {
    $ws.Cells.Item($r, 6).Value = $elisFoto
    $ws.Cells.Item($r, 7).Value = $elisNome
}

# --- Drop the "AA TESTE" scratch rows (406-409) ----------------------------
$ws.Range("A406:G409").Delete()
